$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Cells whose style changes from blank placeholder to numeric/percent ---
$ws.Range("C20").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("C20").Value = 2
$ws.Range("D22").NumberFormat = $ws.Range("F22").NumberFormat
$ws.Range("D22").Value = 3
$ws.Range("E22").NumberFormat = $ws.Range("H22").NumberFormat
$ws.Range("E22").Value = -100
$ws.Range("C23").NumberFormat = $ws.Range("D23").NumberFormat
$ws.Range("C23").Value = 1
$ws.Range("C28").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("C28").Value = 1

# --- Plain value updates ---
# Row 15
$ws.Range("M15").Value = 150

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = -21.428571428571
$ws.Range("L16").Value = -28.260869565217
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = -82.446808510638

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = -21.951219512195
$ws.Range("L17").Value = -13.513513513513
$ws.Range("M17").Value = -25.581395348837
$ws.Range("N17").Value = -56.756756756756

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = -11.764705882352
$ws.Range("L18").Value = -26.829268292682
$ws.Range("M18").Value = -6.25
$ws.Range("N18").Value = -78.873239436619

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = -1.666666666666
$ws.Range("I19").Value = 200
$ws.Range("J19").Value = 222
$ws.Range("K19").Value = -9.909909909909
$ws.Range("L19").Value = 3.626943005181
$ws.Range("M19").Value = 10.497237569060
$ws.Range("N19").Value = -21.259842519685

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -62.5
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 6.25
$ws.Range("L20").Value = -5.555555555555
$ws.Range("M20").Value = 13.333333333333
$ws.Range("N20").Value = -90.395480225988

# Row 21
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 93
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -13.084112149532
$ws.Range("I21").Value = 317
$ws.Range("J21").Value = 357
$ws.Range("K21").Value = -11.204481792717
$ws.Range("L21").Value = -6.764705882352
$ws.Range("M21").Value = 7.094594594594
$ws.Range("N21").Value = -62.440758293838

# Row 22
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -15.384615384615
$ws.Range("L22").Value = 22.222222222222

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = -20
$ws.Range("L23").Value = -42.857142857142
$ws.Range("M23").Value = -25

# Row 24
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 42
$ws.Range("H24").Value = -22.222222222222
$ws.Range("I24").Value = 176
$ws.Range("J24").Value = 199
$ws.Range("K24").Value = -11.557788944723
$ws.Range("L24").Value = -6.878306878306
$ws.Range("M24").Value = -26.050420168067

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 80
$ws.Range("J25").Value = 118
$ws.Range("K25").Value = -32.203389830508

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 83.333333333333
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 13.636363636363
$ws.Range("I26").Value = 101
$ws.Range("J26").Value = 105
$ws.Range("K26").Value = -3.809523809523
$ws.Range("L26").Value = 32.894736842105
$ws.Range("M26").Value = 3.061224489795

# Row 27
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = 233.333333333333
$ws.Range("L27").Value = 150

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 19
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = -9.523809523809
$ws.Range("L28").Value = 18.75

# Row 29
$ws.Range("N29").Value = -80

# Row 30
$ws.Range("N30").Value = -80

# Row 31
$ws.Range("J31").Value = 10
$ws.Range("K31").Value = -90
$ws.Range("L31").Value = -83.333333333333
